{"js": "// Office.js (Word JavaScript API) script.\n// Implements the two edits described by the diff:\n//   1. Merge the two runs of the \"Git Remote\" description paragraph into a\n//      single run with the combined text.\n//   2. Delete the leftover \"Git add\" block (4 blank paragraphs, the\n//      \"Git add\" heading paragraph, and the \"-&gt;Git add ...\" paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Merge the \"Git remote\" description runs into a single run -------\nconst remoteNeedle = \"->Git remote is a command used in Git that\";\nlet remoteParagraph = null;\nfor (const p of paragraphs.items) {\n    if (p.text.indexOf(remoteNeedle) !== -1) {\n        remoteParagraph = p;\n        break;\n    }\n}\n\nif (remoteParagraph) {\n    const fullText =\n        \"->Git remote is a command used in Git that can help you manage connections to remote repositories.\";\n    const range = remoteParagraph.getRange(\"Whole\");\n    range.insertText(fullText, \"Replace\");\n    await context.sync();\n}\n\n// --- 2. Delete the obsolete \"Git add\" paragraphs --------------------------\n// Re-load paragraphs since the collection/text changed above.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text,items/alignment\");\nawait context.sync();\n\nconst items = paragraphs2.items;\n\n// Find the \"Git add\" heading paragraph and the paragraph that starts the\n// \"->Git add\" description right after it.\nlet addHeadingIndex = -1;\nlet addBodyIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n    const t = items[i].text.trim();\n    if (t === \"Git add\") {\n        addHeadingIndex = i;\n    } else if (addHeadingIndex !== -1 && t.indexOf(\"->Git add\") === 0) {\n        addBodyIndex = i;\n        break;\n    }\n}\n\nif (addHeadingIndex !== -1 && addBodyIndex !== -1) {\n    // Walk backwards from the heading paragraph, collecting the run of\n    // blank, identically-(Justified-)aligned paragraphs immediately\n    // preceding it (there are 4 of them in this document). Stop as soon as\n    // either the text is non-empty or the alignment differs, so earlier,\n    // unrelated blank paragraphs (e.g. around the screenshot image) are\n    // left untouched.\n    const headingAlignment = items[addHeadingIndex].alignment;\n    let firstEmptyIndex = addHeadingIndex;\n    while (\n        firstEmptyIndex - 1 >= 0 &&\n        items[firstEmptyIndex - 1].text.trim() === \"\" &&\n        items[firstEmptyIndex - 1].alignment === headingAlignment\n    ) {\n        firstEmptyIndex--;\n    }\n\n    // Delete from the last paragraph to the first so indices stay valid.\n    for (let i = addBodyIndex; i >= firstEmptyIndex; i--) {\n        items[i].delete();\n    }\n    await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Implements the two edits described by the diff:\n#   1. Merge the two runs of the \"Git Remote\" description paragraph into a\n#      single run with the combined text.\n#   2. Delete the leftover \"Git add\" block (4 blank paragraphs, the\n#      \"Git add\" heading paragraph, and the \"->Git add ...\" paragraph).\n\n$d = $word.ActiveDocument\n\n# --- 1. Merge the \"Git remote\" description runs into a single run --------\n$remoteNeedle = \"Git remote is a command used in Git that\"\n$remoteParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains($remoteNeedle)) {\n        $remoteParaIndex = $i\n        break\n    }\n}\n\nif ($remoteParaIndex -ne -1) {\n    $p = $d.Paragraphs.Item($remoteParaIndex)\n    $fullText = \"->Git remote is a command used in Git that can help you manage connections to remote repositories.\"\n\n    # Capture the formatting of the first run so we can re-apply it once the\n    # paragraph's runs are collapsed into one.\n    $fontName = $p.Range.Font.Name\n    $fontSize = $p.Range.Font.Size\n\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End\n\n    # Delete the whole paragraph's text (but not its paragraph mark), then\n    # insert the combined text and restore formatting on it.\n    $textRange = $d.Range($pStart, $pEnd - 1)\n    $textRange.Delete()\n\n    $p2 = $d.Paragraphs.Item($remoteParaIndex)\n    $insertRange = $p2.Range\n    $insertRange.InsertAfter($fullText)\n\n    $p3 = $d.Paragraphs.Item($remoteParaIndex)\n    $newRange = $p3.Range\n    $newRange.Font.Name = $fontName\n    $newRange.Font.Size = $fontSize\n    $newRange.Font.SizeBi = $fontSize\n}\n\n# --- 2. Delete the obsolete \"Git add\" paragraphs --------------------------\n$addHeadingIndex = -1\n$addBodyIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Git add\") {\n        $addHeadingIndex = $i\n    } elseif ($addHeadingIndex -ne -1 -and $t.StartsWith(\"->Git add\")) {\n        $addBodyIndex = $i\n        break\n    }\n}\n\nif ($addHeadingIndex -ne -1 -and $addBodyIndex -ne -1) {\n    # Walk backwards from the heading paragraph, collecting the run of\n    # blank, identically-aligned paragraphs immediately preceding it (there\n    # are 4 of them in this document). Stop as soon as either the text is\n    # non-empty or the alignment differs, so earlier, unrelated blank\n    # paragraphs (e.g. around the screenshot image) are left untouched.\n    $headingAlignment = $d.Paragraphs.Item($addHeadingIndex).Alignment\n    $firstEmptyIndex = $addHeadingIndex\n    while ($firstEmptyIndex - 1 -ge 1) {\n        $prev = $d.Paragraphs.Item($firstEmptyIndex - 1)\n        if ($prev.Range.Text.Trim() -eq \"\" -and $prev.Alignment -eq $headingAlignment) {\n            $firstEmptyIndex = $firstEmptyIndex - 1\n        } else {\n            break\n        }\n    }\n\n    $pFirst = $d.Paragraphs.Item($firstEmptyIndex)\n    $pLast = $d.Paragraphs.Item($addBodyIndex)\n    $delRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)\n    $delRange.Delete()\n}\n"}
